# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "GENESIS TORRES RICO" worker block (periods 2304..2312) is moved to the
# top of the data table (rows 16-24, newest period first) with an updated
# "Salario Basico" (column G) of 1,160,000 (was 1,500,000), and the F16 (Valor
# Mora for period 2312) value is set to 38000.
# The previously-existing workers "KEIVER BOSSIO BALLESTEROS" (period 2212)
# and "CINDY MILENA LARA ESPITALETA" (period 2303) are moved down to become
# the last two rows (25-26) of the table, retaining their original values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout for the data table (rows 16-26):
#  B = Tipo Doc Trabajador, C = N Doc Trabajador, D = Nombre Trabajador,
#  E = Periodo Mora, F = Valor Mora, G = Salario Basico

$rows = @(
    @{ Row = 16; Doc = "1044928283"; Nombre = "GENESIS TORRES RICO"; Periodo = "2312"; Mora = 38000;  Salario = 1160000 },
    @{ Row = 17; Doc = "1044928283"; Nombre = "GENESIS TORRES RICO"; Periodo = "2311"; Mora = 60000;  Salario = 1160000 },
    @{ Row = 18; Doc = "1044928283"; Nombre = "GENESIS TORRES RICO"; Periodo = "2310"; Mora = 60000;  Salario = 1160000 },
    @{ Row = 19; Doc = "1044928283"; Nombre = "GENESIS TORRES RICO"; Periodo = "2309"; Mora = 60000;  Salario = 1160000 },
    @{ Row = 20; Doc = "1044928283"; Nombre = "GENESIS TORRES RICO"; Periodo = "2308"; Mora = 60000;  Salario = 1160000 },
    @{ Row = 21; Doc = "1044928283"; Nombre = "GENESIS TORRES RICO"; Periodo = "2307"; Mora = 60000;  Salario = 1160000 },
    @{ Row = 22; Doc = "1044928283"; Nombre = "GENESIS TORRES RICO"; Periodo = "2306"; Mora = 60000;  Salario = 1160000 },
    @{ Row = 23; Doc = "1044928283"; Nombre = "GENESIS TORRES RICO"; Periodo = "2305"; Mora = 60000;  Salario = 1160000 },
    @{ Row = 24; Doc = "1044928283"; Nombre = "GENESIS TORRES RICO"; Periodo = "2304"; Mora = 60000;  Salario = 1160000 },
    @{ Row = 25; Doc = "1002319881"; Nombre = "KEIVER BOSSIO BALLESTEROS";    Periodo = "2212"; Mora = 34666; Salario = 1000000 },
    @{ Row = 26; Doc = "1050952836"; Nombre = "CINDY MILENA LARA ESPITALETA"; Periodo = "2303"; Mora = 16000; Salario = 1500000 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("B$n").Value = "CC"
    $ws.Range("C$n").Value = $r.Doc
    $ws.Range("D$n").Value = $r.Nombre
    $ws.Range("E$n").Value = $r.Periodo
    $ws.Range("F$n").Value = $r.Mora
    $ws.Range("G$n").Value = $r.Salario
}
